$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model results")

$ws.Range("A19").Value = "Rhopalosiphum rufiabdominalis US Weston"
$ws.Range("H19").Value = 0
$ws.Range("K19").Value = "Including overwintering; DDE model run for 15 years for historical period; no future time series data b/c extinct"

$ws.Range("K23").Value = "Including overwintering; extinct in model"
$ws.Range("K24").Value = "Including overwintering; extinct in model"
$ws.Range("K27").Value = "Including overwintering; extinct in model"
$ws.Range("K28").Value = "Including overwintering; extinct in model"
$ws.Range("K29").Value = "Including overwintering; extinct in model"

[void]$ws.Range("K23").Select()
